# pptx: Include image title in description
#
# The image title (i.e. `![alt text](link "title")`) was previously ignored
# when writing to PowerPoint's description (alt text) of a picture, which
# only carried the link/filename (e.g. "lalune.jpg"). This script folds the
# title back in, the same way pandoc now does: "fig:  lalune.jpg" (the
# title "fig:" followed by the link).
#
# Only the two slides whose picture is a standalone figure (i.e. was given
# a title) are affected; the lone full-bleed picture that is the sole shape
# on the first slide had no title, so its description is left untouched.

$p = $ppt.ActivePresentation

function Set-FigureDescription($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.Type -eq 13 -and $shape.AlternativeText -eq "lalune.jpg") {
            $shape.AlternativeText = "fig:  lalune.jpg"
        }
    }
}

Set-FigureDescription 2
Set-FigureDescription 3
